$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.930.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.35%  "

$ws.Range("D3").Value = "'1.735.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.71%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'309.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.02%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4982"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("D8").Value = "'0.3519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").Value = "'42.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("D10").Value = "'0.07237"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("D11").Value = "'1.055"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.73%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "'19.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "

$ws.Range("D14").Value = "'5.938"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "'1.739.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("D16").Value = "'6.852"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.68%  "

$ws.Range("D17").Value = "'86.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.69%  "

$ws.Range("D18").Value = "'0.00001034"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.83%  "

$ws.Range("D19").Value = "'0.06406"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "'16.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").Value = "'5.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").Value = "'27.023.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("E24").Value = "  -1.73%  "

$ws.Range("D25").Value = "'2.064"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.52%  "

$ws.Range("D26").Value = "'153.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.43%  "

$ws.Range("D27").Value = "'19.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").Value = "'1.939.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("D29").Value = "'2.100"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.66%  "

$ws.Range("D30").Value = "'120.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.79%  "

$ws.Range("D31").Value = "'1.057"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "'0.09415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").Value = "'3.575"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").Value = "'5.381"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").Value = "'0.05918"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("D36").Value = "'0.02186"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.50%  "

$ws.Range("D37").Value = "'10.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.08%  "

$ws.Range("D38").Value = "'1.427"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.740"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.53%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.1983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("D41").Value = "'0.6016"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").Value = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'1.109"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'7.424"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.19%  "

$ws.Range("D45").Value = "'12.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").Value = "'3.567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.63%  "

$ws.Range("D47").Value = "'0.5629"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.09%  "

$ws.Range("D48").Value = "'118.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("D49").Value = "'1.843"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.78%  "

$ws.Range("D50").Value = "'0.06662"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("D51").Value = "'1.095"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.28%  "
